# "A Big Performance Fix"
#
# The card table's per-row "default asset" columns (K:Q = Address of Image,
# Address of Idle Gif, Address of Run Gif, Address of Attack Gif,
# Address of Get Damage Gif, Address Of Death Gif, Target Society) were
# left as the literal placeholder "Null" for every card row except the
# template row (row 2). Populate every card row (3-51) with the same
# default asset references used by row 2, removing the "Null" placeholder
# entirely (the game was presumably doing a slow per-row fallback/lookup
# whenever it hit "Null").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$imagePath  = "./res/Characters/generals/general_f1.png"
$idleGif    = "./res/gifs/f1_altgeneral/idle_t.gif"
$runGif     = "./res/gifs/f1_altgeneral/run_t.gif"
$attackGif  = "./res/gifs/f1_altgeneral/attack_t.gif"
$hitGif     = "./res/gifs/f1_altgeneral/hit_t.gif"
$deathGif   = "./res/gifs/f1_altgeneral/death_t.gif"
$society    = "FRIENDLY"

for ($r = 3; $r -le 51; $r++) {
    $ws.Cells.Item($r, 11).Value = $imagePath   # K: Address of Image
    $ws.Cells.Item($r, 12).Value = $idleGif      # L: Address of Idle Gif
    $ws.Cells.Item($r, 13).Value = $runGif       # M: Address of Run Gif
    $ws.Cells.Item($r, 14).Value = $attackGif    # N: Address of Attack Gif
    $ws.Cells.Item($r, 15).Value = $hitGif       # O: Address of Get Damage Gif
    $ws.Cells.Item($r, 16).Value = $deathGif     # P: Address Of Death Gif
    $ws.Cells.Item($r, 17).Value = $society      # Q: Target Society
}

# Widen column B (Card Description) substantially so the long Farsi
# description strings are readable, and give the new scroll position /
# selection used while reviewing the fix.
$ws.Columns.Item(2).ColumnWidth = 89.6

$ws.Range("P54").Select()
